$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 26,14
$data[0,0] = 1
$data[0,1] = 'bedrooms'
$data[0,2] = 'bedrooms'
$data[0,3] = 'target'
$data[0,5] = 'j'
$data[0,6] = 'stimuli/img_yteqw.png'
$data[0,7] = 66.83783783783784
$data[0,8] = 43.78378378378378
$data[0,9] = 55.31081081081081
$data[0,10] = 37
$data[0,11] = 4
$data[0,12] = 4
$data[0,13] = 4
$data[1,0] = 2
$data[1,1] = 'bedrooms'
$data[1,2] = 'kitchens'
$data[1,3] = 'distractor'
$data[1,5] = 'f'
$data[1,6] = 'stimuli/img_411xa.png'
$data[1,7] = 51.03030303030303
$data[1,8] = 28.93939393939394
$data[1,9] = 39.98484848484848
$data[1,10] = 33
$data[1,11] = 2
$data[1,12] = 2
$data[1,13] = 2
$data[2,0] = 3
$data[2,1] = 'bedrooms'
$data[2,2] = 'living_rooms'
$data[2,3] = 'distractor'
$data[2,5] = 'f'
$data[2,6] = 'stimuli/img_gka64.png'
$data[2,7] = 19.23809523809524
$data[2,8] = 20.02380952380953
$data[2,9] = 19.63095238095238
$data[2,10] = 42
$data[2,11] = 1
$data[2,12] = 1
$data[2,13] = 1
$data[3,0] = 4
$data[3,1] = 'bedrooms'
$data[3,2] = 'bedrooms'
$data[3,3] = 'target'
$data[3,5] = 'j'
$data[3,6] = 'stimuli/img_72fmj.png'
$data[3,7] = 53.87179487179487
$data[3,8] = 36.02564102564103
$data[3,9] = 44.94871794871795
$data[3,10] = 39
$data[3,11] = 3
$data[3,12] = 3
$data[3,13] = 3
$data[4,0] = 5
$data[4,1] = 'bedrooms'
$data[4,2] = 'kitchens'
$data[4,3] = 'distractor'
$data[4,5] = 'f'
$data[4,6] = 'stimuli/img_eppte.png'
$data[4,7] = 78.42424242424242
$data[4,8] = 57.03030303030303
$data[4,9] = 67.72727272727272
$data[4,10] = 33
$data[4,11] = 7
$data[4,12] = 7
$data[4,13] = 7
$data[5,0] = 6
$data[5,1] = 'bedrooms'
$data[5,2] = 'bedrooms'
$data[5,3] = 'target'
$data[5,5] = 'j'
$data[5,6] = 'stimuli/img_9pfbj.png'
$data[5,7] = 91.27272727272727
$data[5,8] = 80.0909090909091
$data[5,9] = 85.68181818181819
$data[5,10] = 33
$data[5,11] = 10
$data[5,12] = 10
$data[5,13] = 10
$data[6,0] = 7
$data[6,1] = 'bedrooms'
$data[6,2] = 'bedrooms'
$data[6,3] = 'target'
$data[6,5] = 'j'
$data[6,6] = 'stimuli/img_kzg3h.png'
$data[6,7] = 77.02777777777777
$data[6,8] = 56.22222222222222
$data[6,9] = 66.625
$data[6,10] = 36
$data[6,11] = 7
$data[6,12] = 7
$data[6,13] = 7
$data[7,0] = 8
$data[7,1] = 'bedrooms'
$data[7,2] = 'bedrooms'
$data[7,3] = 'target'
$data[7,5] = 'j'
$data[7,6] = 'stimuli/img_ose78.png'
$data[7,7] = 80.19444444444444
$data[7,8] = 60.25
$data[7,9] = 70.22222222222223
$data[7,10] = 36
$data[7,11] = 8
$data[7,12] = 7
$data[7,13] = 7
$data[8,0] = 9
$data[8,1] = 'bedrooms'
$data[8,2] = 'bedrooms'
$data[8,3] = 'target'
$data[8,5] = 'j'
$data[8,6] = 'stimuli/img_3bxjb.png'
$data[8,7] = 87.28571428571429
$data[8,8] = 72.65714285714286
$data[8,9] = 79.97142857142858
$data[8,10] = 35
$data[8,11] = 10
$data[8,12] = 10
$data[8,13] = 10
$data[9,0] = 10
$data[9,1] = 'bedrooms'
$data[9,2] = 'bedrooms'
$data[9,3] = 'target'
$data[9,5] = 'j'
$data[9,6] = 'stimuli/img_f4jxo.png'
$data[9,7] = 82.91666666666667
$data[9,8] = 65.52777777777777
$data[9,9] = 74.22222222222223
$data[9,10] = 36
$data[9,11] = 8
$data[9,12] = 8
$data[9,13] = 8
$data[10,0] = 11
$data[10,1] = 'bedrooms'
$data[10,2] = 'bedrooms'
$data[10,3] = 'target'
$data[10,5] = 'j'
$data[10,6] = 'stimuli/img_jivhq.png'
$data[10,7] = 37
$data[10,8] = 22.26530612244898
$data[10,9] = 29.63265306122449
$data[10,10] = 49
$data[10,11] = 2
$data[10,12] = 2
$data[10,13] = 2
$data[11,0] = 12
$data[11,1] = 'bedrooms'
$data[11,2] = 'kitchens'
$data[11,3] = 'distractor'
$data[11,5] = 'f'
$data[11,6] = 'stimuli/img_xguy9.png'
$data[11,7] = 78.21621621621621
$data[11,8] = 57.24324324324324
$data[11,9] = 67.72972972972973
$data[11,10] = 37
$data[11,11] = 7
$data[11,12] = 7
$data[11,13] = 7
$data[12,0] = 13
$data[12,1] = 'bedrooms'
$data[12,2] = 'bedrooms'
$data[12,3] = 'target'
$data[12,5] = 'j'
$data[12,6] = 'stimuli/img_t4hvr.png'
$data[12,7] = 61.69230769230769
$data[12,8] = 39.76923076923077
$data[12,9] = 50.73076923076923
$data[12,10] = 39
$data[12,11] = 3
$data[12,12] = 3
$data[12,13] = 3
$data[13,0] = 14
$data[13,1] = 'bedrooms'
$data[13,2] = 'bedrooms'
$data[13,3] = 'target'
$data[13,5] = 'j'
$data[13,6] = 'stimuli/img_cmyvx.png'
$data[13,7] = 64.25
$data[13,8] = 40.09375
$data[13,9] = 52.171875
$data[13,10] = 32
$data[13,11] = 4
$data[13,12] = 4
$data[13,13] = 4
$data[14,0] = 15
$data[14,1] = 'bedrooms'
$data[14,2] = 'kitchens'
$data[14,3] = 'distractor'
$data[14,5] = 'f'
$data[14,6] = 'stimuli/img_q1ynd.png'
$data[14,7] = 70.05714285714286
$data[14,8] = 47.31428571428572
$data[14,9] = 58.68571428571429
$data[14,10] = 35
$data[14,11] = 5
$data[14,12] = 5
$data[14,13] = 5
$data[15,0] = 16
$data[15,1] = 'bedrooms'
$data[15,2] = 'bedrooms'
$data[15,3] = 'target'
$data[15,5] = 'j'
$data[15,6] = 'stimuli/img_ic3os.png'
$data[15,7] = 84.79069767441861
$data[15,8] = 66.16279069767442
$data[15,9] = 75.47674418604652
$data[15,10] = 43
$data[15,11] = 9
$data[15,12] = 9
$data[15,13] = 9
$data[16,0] = 17
$data[16,1] = 'bedrooms'
$data[16,2] = 'bedrooms'
$data[16,3] = 'target'
$data[16,5] = 'j'
$data[16,6] = 'stimuli/img_z3yzz.png'
$data[16,7] = 71.71052631578948
$data[16,8] = 49.81578947368421
$data[16,9] = 60.76315789473685
$data[16,10] = 38
$data[16,11] = 5
$data[16,12] = 5
$data[16,13] = 5
$data[17,0] = 18
$data[17,1] = 'bedrooms'
$data[17,2] = 'bedrooms'
$data[17,3] = 'target'
$data[17,5] = 'j'
$data[17,6] = 'stimuli/img_aweye.png'
$data[17,7] = 53.42105263157895
$data[17,8] = 31.84210526315789
$data[17,9] = 42.63157894736842
$data[17,10] = 38
$data[17,11] = 2
$data[17,12] = 2
$data[17,13] = 2
$data[18,0] = 19
$data[18,1] = 'bedrooms'
$data[18,2] = 'bedrooms'
$data[18,3] = 'target'
$data[18,5] = 'j'
$data[18,6] = 'stimuli/img_1vq1v.png'
$data[18,7] = 69.42857142857143
$data[18,8] = 46.59523809523809
$data[18,9] = 58.01190476190476
$data[18,10] = 42
$data[18,11] = 5
$data[18,12] = 5
$data[18,13] = 5
$data[19,0] = 20
$data[19,1] = 'bedrooms'
$data[19,2] = 'bedrooms'
$data[19,3] = 'target'
$data[19,5] = 'j'
$data[19,6] = 'stimuli/img_anzgh.png'
$data[19,7] = 75.10526315789474
$data[19,8] = 55.76315789473684
$data[19,9] = 65.4342105263158
$data[19,10] = 38
$data[19,11] = 6
$data[19,12] = 6
$data[19,13] = 6
$data[20,0] = 21
$data[20,1] = 'bedrooms'
$data[20,2] = 'bedrooms'
$data[20,3] = 'target'
$data[20,5] = 'j'
$data[20,6] = 'stimuli/img_gbypq.png'
$data[20,7] = 76.275
$data[20,8] = 51.925
$data[20,9] = 64.1
$data[20,10] = 40
$data[20,11] = 6
$data[20,12] = 6
$data[20,13] = 6
$data[21,0] = 22
$data[21,1] = 'bedrooms'
$data[21,2] = 'bedrooms'
$data[21,3] = 'target'
$data[21,5] = 'j'
$data[21,6] = 'stimuli/img_fnu4h.png'
$data[21,7] = 85.87179487179488
$data[21,8] = 70.71794871794872
$data[21,9] = 78.2948717948718
$data[21,10] = 39
$data[21,11] = 9
$data[21,12] = 9
$data[21,13] = 9
$data[22,0] = 23
$data[22,1] = 'bedrooms'
$data[22,2] = 'bedrooms'
$data[22,3] = 'target'
$data[22,5] = 'j'
$data[22,6] = 'stimuli/img_2pnl2.png'
$data[22,7] = 6.621621621621622
$data[22,8] = 7.135135135135135
$data[22,9] = 6.878378378378379
$data[22,10] = 37
$data[22,11] = 1
$data[22,12] = 1
$data[22,13] = 1
$data[23,0] = 24
$data[23,1] = 'bedrooms'
$data[23,2] = 'bedrooms'
$data[23,3] = 'target'
$data[23,5] = 'j'
$data[23,6] = 'stimuli/img_cgdyc.png'
$data[23,7] = 32.93023255813954
$data[23,8] = 14.04651162790698
$data[23,9] = 23.48837209302326
$data[23,10] = 43
$data[23,11] = 1
$data[23,12] = 1
$data[23,13] = 1
$data[24,0] = 25
$data[24,1] = 'bedrooms'
$data[24,2] = 'bedrooms'
$data[24,3] = 'target'
$data[24,5] = 'j'
$data[24,6] = 'stimuli/img_juob3.png'
$data[24,7] = 79.92105263157895
$data[24,8] = 59.78947368421053
$data[24,9] = 69.85526315789474
$data[24,10] = 38
$data[24,11] = 7
$data[24,12] = 7
$data[24,13] = 7
$data[25,0] = 26
$data[25,1] = 'bedrooms'
$data[25,2] = 'living_rooms'
$data[25,3] = 'distractor'
$data[25,5] = 'f'
$data[25,6] = 'stimuli/img_16kib.png'
$data[25,7] = 80.97727272727273
$data[25,8] = 61.11363636363637
$data[25,9] = 71.04545454545455
$data[25,10] = 44
$data[25,11] = 1
$data[25,12] = 1
$data[25,13] = 1

$ws.Range("F2:S27").Value = $data
